# Edit script: Tarea 1 entry + "Nelmer P." run cleanup + namespace cleanup.
#
# Strategy: the target XML keeps each paragraph's original w14:paraId /
# w:rsidR / w:rsidRDefault attributes untouched and just changes the runs
# inside it, so rather than typing text through Range.Text (which drops the
# explicit w:rFonts on the newly created run) we rebuild each paragraph with
# Range.InsertXML, re-specifying the same paraId/rsid attributes the
# paragraph already had. InsertXML replaces the whole enclosing <w:p>, which
# is why we carry those attributes along explicitly.

$d = $word.ActiveDocument

function New-WordPackageXml([string]$bodyInner) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$arialRpr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr>'

# ---------------------------------------------------------------------
# 1) "Nelmer" + " P." (split across two runs, wrapped in proofErr spell
#    check markers) -> single run "Nelmer P." with no proofErr markers.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Nelmer") | Out-Null
$nelmerStart = $rng.Start
$rng2 = $d.Content
$rng2.Find.Execute("Camilo") | Out-Null
$camiloStart = $rng2.Start
# The "Nelmer P." paragraph runs from the start of "Nelmer" up to (but not
# including) the start of the following "Camilo" paragraph, i.e. it also
# swallows its own paragraph mark.
$nelmerParaRange = $d.Range($nelmerStart, $camiloStart)

$nelmerP = '<w:p w14:paraId="56200208" w14:textId="77777777" w:rsidR="00194D8A" w:rsidRDefault="00194D8A" w:rsidP="00194D8A"><w:pPr>' + $arialRpr + '</w:pPr><w:r>' + $arialRpr + '<w:t>Nelmer P.</w:t></w:r></w:p>'
$nelmerParaRange.InsertXML((New-WordPackageXml $nelmerP)) | Out-Null

# ---------------------------------------------------------------------
# 2) Fill in the "Tarea 1" row (currently four empty cells).
# ---------------------------------------------------------------------
$t = $d.Tables.Item(1)

# Column 1: Numero -> "1"
$c1 = $t.Cell(3, 1).Range
$p1 = '<w:p w14:paraId="7D7E0487" w14:textId="77777777" w:rsidR="00194D8A" w:rsidRDefault="00194D8A"><w:pPr>' + $arialRpr + '</w:pPr><w:r>' + $arialRpr + '<w:t>1</w:t></w:r></w:p>'
$c1.InsertXML((New-WordPackageXml $p1)) | Out-Null

# Column 2: Fecha -> "6/05/2025"
$c2 = $t.Cell(3, 2).Range
$p2 = '<w:p w14:paraId="6680B9A9" w14:textId="77777777" w:rsidR="00194D8A" w:rsidRDefault="00194D8A"><w:pPr>' + $arialRpr + '</w:pPr><w:r>' + $arialRpr + '<w:t>6/05/2025</w:t></w:r></w:p>'
$c2.InsertXML((New-WordPackageXml $p2)) | Out-Null

# Column 3: Descripcion -> four runs, alternating Arial / default font.
$c3 = $t.Cell(3, 3).Range
$descRuns = ''
$descRuns += '<w:r>' + $arialRpr + '<w:t>Se completo la tarea 1, además de hacer los diagramas de procesos, d</w:t></w:r>'
$descRuns += '<w:r><w:t>e caso de uso</w:t></w:r>'
$descRuns += '<w:r>' + $arialRpr + '<w:t xml:space="preserve"> y de </w:t></w:r>'
$descRuns += '<w:r><w:t>diagrama de secuencia</w:t></w:r>'
$p3 = '<w:p w14:paraId="735F63BB" w14:textId="77777777" w:rsidR="00194D8A" w:rsidRDefault="00194D8A"><w:pPr>' + $arialRpr + '</w:pPr>' + $descRuns + '</w:p>'
$c3.InsertXML((New-WordPackageXml $p3)) | Out-Null

# Column 4: Autor/es -> two paragraphs, "Marlon D." and "Edinson P."
$c4 = $t.Cell(3, 4).Range
$p4a = '<w:p w14:paraId="2D81E6E2" w14:textId="77777777" w:rsidR="00194D8A" w:rsidRDefault="00194D8A"><w:pPr>' + $arialRpr + '</w:pPr><w:r>' + $arialRpr + '<w:t>Marlon D.</w:t></w:r></w:p>'
$p4b = '<w:p><w:pPr>' + $arialRpr + '</w:pPr><w:r>' + $arialRpr + '<w:t>Edinson P.</w:t></w:r></w:p>'
$c4.InsertXML((New-WordPackageXml ($p4a + $p4b))) | Out-Null
